# Refresh the cryptos price/volume table (columns D "Price" and E
# "Volume(1h)") for rows 2-51 with the latest scraped values.
#
# D-column prices are stored as TEXT (e.g. "27.783.11", "1.030", "5.000")
# rather than numbers, so trailing zeros and thousands-dot groupings are
# preserved exactly. Plainly assigning a numeric-looking string to
# .Value lets Excel auto-convert it to a real number (losing trailing
# zeros, e.g. "6.200" -> 6.2), so each D-cell is forced to Text via a
# temporary "@" NumberFormat before the assignment, then restored to the
# workbook's default "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.734.26'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.779.30'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.24%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.29%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5116'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3788'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07776'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.15'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.082'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.200'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.09'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.776.62'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.155'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.21'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.82%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001070'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06554'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.85%  '
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.98'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.909'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.795.95'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.232'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.28'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.20'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.986.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.349'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.85'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.92%  '
$ws.Range('E31').Value = '  -1.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.029'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.629'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.465'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07051'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02312'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.691'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2118'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '11.50'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.001'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6069'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.07%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.149'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.74%  '
$ws.Range('E44').Value = '  -5.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.10'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5911'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.706'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '127.45'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.201'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.889'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06799'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.56%  '
